# RegisterMappings.xlsx edit script
# Commit message: "working on verilog comments. Asm code is done."
#
# Semantic changes:
#  1. Remove the "Seven Seg Upper (7:4) Digits value address" note from
#     C17 ($t7 row) - that register mapping note is being cleared.
#  2. Add a new note "Decimal Point Value register" to C24 ($s6 row).
#  3. Finish the BusBlaster clock-divider table (row 39, the "10hz" row)
#     with the same formula pattern used by the 5hz/1hz rows above it.
#  4. Add a small scratch calculation block below the table (rows 42-45)
#     working out a clock-divider value (50,000,000 Hz / 100 -> 500,000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the stale "Seven Seg Upper" note in C17.
$ws.Range("C17").ClearContents()

# 2. New note for the Decimal Point Value register in C24.
$ws.Range("C24").Value = "Decimal Point Value register"

# 3. Complete row 39 (10hz row) the same way rows 38/39 above were built:
#    D = D37 / 10, E = $D$39 * E36, F:M = $D$39 * (F36:M36) as a fill-right.
$ws.Range("D39").Formula = '=D37/10'
$ws.Range("E39").Formula = '=$D$39*E36'
$ws.Range("F39:M39").Formula = '=$D$39*F36'

# 4. New scratch block under the table.
$ws.Range("C42").Value = 100
$ws.Range("C43").Value = 50000000
$ws.Range("C43").NumberFormat = "0.00E+00"
$ws.Range("C44").Formula = '=C42/C43'
$ws.Range("C44").NumberFormat = "General"
$ws.Range("F44").Formula = '=1/C44'
$ws.Range("J45").Value = 1
$ws.Range("K45").Value = 1
$ws.Range("L45").Value = 1

# Restore the selection to where the author's cursor ended up.
$ws.Range("F25").Select()
